# [MOD] Strengthen Test final Assertion
#
# The original paragraph reads "Sample for Test" as a single run.
# The target content is the same paragraph split into three runs -
# "Find This ", "For", " Test" - with a pair of <w:proofErr> gramStart /
# gramEnd markers (as Word's grammar checker would insert) bracketing
# the word "For".
#
# We rebuild the paragraph's WordprocessingML directly via Range.InsertXML
# so the exact run/proofErr structure is produced, preserving the
# paragraph's existing identity attributes.

$d = $word.ActiveDocument
$p = $d.Paragraphs(1)
$r = $p.Range

$paraXml = '<w:p ' +
    'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
    'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
    'w14:paraId="19EEA3C2" w14:textId="04FE5BA6" ' +
    'w:rsidR="00067157" w:rsidRDefault="00365033">' +
        '<w:r><w:t xml:space="preserve">Find This </w:t></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:t>For</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:t xml:space="preserve"> Test</w:t></w:r>' +
    '</w:p>'

$r.InsertXML($paraXml)
